$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns stay as text so values like "1.0000" or "0.08512"
# are not auto-converted to numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.009.94"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.50%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.659.68"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.36%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.52%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.68"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.78%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.0000"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.68%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3900"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3867"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.48%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.39"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.369"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.41%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9997"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.74%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08512"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.96"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.43%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.222"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.052"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +8.79%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001314"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.68%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.658.89"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.47%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.51"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06991"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.36%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.98"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.59%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.981"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.61%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.67"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.54%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.020.50"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.52%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.492"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.54%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.099"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +8.62%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.25"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "154.08"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.92%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "140.30"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.31%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.322"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.053"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.26%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.489"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.839.08"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.049"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +10.22%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08151"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.96%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.03013"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.07%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +9.41%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.776"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.90%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2711"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.62%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09160"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.74"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.99%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7575"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.422"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.71%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.59"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.92%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7031"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.50%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.503"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.13%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.54%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.000"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08307"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.95%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "135.78"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.58%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.242"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.38%  "
